$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pair output refreshed with new TPM-based values:
#  - Sheet1 now has 8 data rows (was 5): both "ECs" and "MuSCs" appear as
#    the Sending cluster (rows 2-5 and 6-9 respectively), each paired with
#    the same 4 Target clusters.
#  - The "Resolving-Mac" target cluster/string is gone; all LR statistics
#    (cols E:T) are recomputed.
#  - Sheet dimension grows from A1:T6 to A1:T9 (3 new rows appended).

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgi1"
$ws.Range("C2").Value = "Adam11"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.000409
$ws.Range("H2").Value = 0.001227
$ws.Range("I2").Value = 0.164521319388576
$ws.Range("J2").Value = 0.228024530756365
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1135006666666667
$ws.Range("N2").Value = 0.340502
$ws.Range("O2").Value = 0.1328118661732309
$ws.Range("P2").Value = 0.1579033525861449
$ws.Range("Q2").Value = 0.00004642177266666667
$ws.Range("R2").Value = 0.000417795954
$ws.Range("S2").Value = 0.02185038345327893
$ws.Range("T2").Value = 0.03600583787831255

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgi1"
$ws.Range("C3").Value = "Adam11"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.000409
$ws.Range("H3").Value = 0.001227
$ws.Range("I3").Value = 0.164521319388576
$ws.Range("J3").Value = 0.228024530756365
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.183029
$ws.Range("N3").Value = 0.549087
$ws.Range("O3").Value = 0.2141698702546853
$ws.Range("P3").Value = 0.2546319204041931
$ws.Range("Q3").Value = 0.00007485886100000001
$ws.Range("R3").Value = 0.000673729749
$ws.Range("S3").Value = 0.03523550962758096
$ws.Range("T3").Value = 0.05806232416575823

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgi1"
$ws.Range("C4").Value = "Adam11"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.000409
$ws.Range("H4").Value = 0.001227
$ws.Range("I4").Value = 0.164521319388576
$ws.Range("J4").Value = 0.228024530756365
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1506706666666666
$ws.Range("N4").Value = 0.452012
$ws.Range("O4").Value = 0.1763060341868607
$ws.Range("P4").Value = 0.209614657796925
$ws.Range("Q4").Value = 0.00006162430266666666
$ws.Range("R4").Value = 0.0005546187239999999
$ws.Range("S4").Value = 0.02900610136058971
$ws.Range("T4").Value = 0.04779728398379984

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lgi1"
$ws.Range("C5").Value = "Adam11"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.000409
$ws.Range("H5").Value = 0.001227
$ws.Range("I5").Value = 0.164521319388576
$ws.Range("J5").Value = 0.228024530756365
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.407397
$ws.Range("N5").Value = 0.814794
$ws.Range("O5").Value = 0.4767122293852232
$ws.Range("P5").Value = 0.3778500692127371
$ws.Range("Q5").Value = 0.000166625373
$ws.Range("R5").Value = 0.0009997522379999999
$ws.Range("S5").Value = 0.07842932494712643
$ws.Range("T5").Value = 0.0861590847284944

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Lgi1"
$ws.Range("C6").Value = "Adam11"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.002077
$ws.Range("H6").Value = 0.004154
$ws.Range("I6").Value = 0.8354786806114239
$ws.Range("J6").Value = 0.771975469243635
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1135006666666667
$ws.Range("N6").Value = 0.340502
$ws.Range("O6").Value = 0.1328118661732309
$ws.Range("P6").Value = 0.1579033525861449
$ws.Range("Q6").Value = 0.0002357408846666666
$ws.Range("R6").Value = 0.001414445308
$ws.Range("S6").Value = 0.1109614827199519
$ws.Range("T6").Value = 0.1218975147078324

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Lgi1"
$ws.Range("C7").Value = "Adam11"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.002077
$ws.Range("H7").Value = 0.004154
$ws.Range("I7").Value = 0.8354786806114239
$ws.Range("J7").Value = 0.771975469243635
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.183029
$ws.Range("N7").Value = 0.549087
$ws.Range("O7").Value = 0.2141698702546853
$ws.Range("P7").Value = 0.2546319204041931
$ws.Range("Q7").Value = 0.0003801512329999999
$ws.Range("R7").Value = 0.002280907398
$ws.Range("S7").Value = 0.1789343606271043
$ws.Range("T7").Value = 0.1965695962384349

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Lgi1"
$ws.Range("C8").Value = "Adam11"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.002077
$ws.Range("H8").Value = 0.004154
$ws.Range("I8").Value = 0.8354786806114239
$ws.Range("J8").Value = 0.771975469243635
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1506706666666666
$ws.Range("N8").Value = 0.452012
$ws.Range("O8").Value = 0.1763060341868607
$ws.Range("P8").Value = 0.209614657796925
$ws.Range("Q8").Value = 0.0003129429746666666
$ws.Range("R8").Value = 0.001877657848
$ws.Range("S8").Value = 0.147299932826271
$ws.Range("T8").Value = 0.1618173738131251

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Lgi1"
$ws.Range("C9").Value = "Adam11"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.002077
$ws.Range("H9").Value = 0.004154
$ws.Range("I9").Value = 0.8354786806114239
$ws.Range("J9").Value = 0.771975469243635
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.407397
$ws.Range("N9").Value = 0.814794
$ws.Range("O9").Value = 0.4767122293852232
$ws.Range("P9").Value = 0.3778500692127371
$ws.Range("Q9").Value = 0.0008461635689999999
$ws.Range("R9").Value = 0.003384654276
$ws.Range("S9").Value = 0.3982829044380967
$ws.Range("T9").Value = 0.2916909844842427

